# Soulreaping weapons cannot be disenchanted (#481)
#
# The "Weapons" sheet lists, per weapon-material, the temper/breakdown
# ingredients and the crafting perk required. Soulreaping weapons were
# missing a row, which meant the patcher treated them as disenchantable.
# Insert a new data row (alphabetically placed between "SkyforgeSteel"
# and "Spectral") describing Soulreaping weapons as craftable (Steel
# Ingot temper/breakdown, Craftsmanship perk) so they can't be broken
# down / disenchanted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weapons")

# Push the existing "Spectral" row (and everything below it) down one,
# opening up row 34 for the new entry.
$ws.Rows.Item(34).Insert()

$ws.Range("A34").Value = "Soulreaping"
$ws.Range("B34").Value = 3
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = 2 / 3
$ws.Range("D34").NumberFormat = "# ?/?"
$ws.Range("G34").Value = "Steel Ingot"
$ws.Range("H34").Value = "Steel Ingot"
$ws.Range("I34").Value = "Craftsmanship"

# Reflect where editing was left off: Weapons becomes the active tab,
# scrolled/selected near the bottom of the newly-extended table.
$null = $ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("A39").Select()
